# Edit script: restructure workbook from 3 sheets (2_50, 3_14, 3_53) for N2 measurements
# into 6 sheets covering both N2 and Ar gas measurements across 4 samples.
$wb = $excel.ActiveWorkbook

# --- Step 1: insert the new "3_14_Ar" sheet before the currently-active sheet (3_53) ---
$sheetAr314 = $wb.Worksheets.Add()
$sheetAr314.Name = "3_14_Ar"

# --- Step 2: rename the pre-existing N2 sheets ---
$wb.Worksheets.Item(1).Name = "2_50_N2"
$wb.Worksheets.Item(2).Name = "3_14_N2"
$wb.Worksheets.Item(4).Name = "3_53_N2"

# --- Step 3: populate "3_14_Ar" with its isotherm data ---
$sheetAr314.Cells.Item(1,1).Value = "p_rels"
$sheetAr314.Cells.Item(1,2).Value = "q_abs"
$data314Ar = @(@(0.30039199357601903,23.485313979899978),@(0.35157831214384816,25.096830687341605),@(0.40083385408149758,26.654272914671349),@(0.45023812996431944,28.328064345279635),@(0.50017180111076964,30.063851833641163),@(0.55050893566871761,31.971497841535967),@(0.60052934519391066,34.086843797606051),@(0.65035155363844099,36.530207219789425),@(0.69967237242768132,39.372755128896202),@(0.74970732790560368,42.820418167589622),@(0.79949086509117706,47.301878156506106),@(0.82436493184885828,50.324152157731945),@(0.85020022918810589,53.963825289124699),@(0.87422796715016959,58.176293881272059),@(0.89901568768798401,63.721176227716434),@(0.9230713909952859,71.004175191356225),@(0.94648540550914972,81.629143906460556),@(0.95785726488730905,89.565133774591999),@(0.96849475642195881,100.10731982164464),@(0.97758804409711664,113.56300189177291),@(0.9860010623627754,134.30653303464081),@(0.9903985579994079,152.57754298777792),@(0.99347609632619271,169.6468955741785))
for ($i = 0; $i -lt $data314Ar.Count; $i++) {
    $r = $i + 2
    $sheetAr314.Cells.Item($r, 1).Value = $data314Ar[$i][0]
    $sheetAr314.Cells.Item($r, 2).Value = $data314Ar[$i][1]
}
$sheetAr314.Range("E14").Select()

# --- Step 4: add "1_223_kerogen_N2" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetKerogenN2 = $wb.Worksheets.Add($null, $lastSheet)
$sheetKerogenN2.Name = "1_223_kerogen_N2"
$sheetKerogenN2.Cells.Item(1,1).Value = "p_rels"
$sheetKerogenN2.Cells.Item(1,2).Value = "q_abs"
$dataKerogenN2 = @(@(0.30119253839666976,43.520344870481068),@(0.35193313464421877,46.189840567168481),@(0.40016151218164148,48.724705505002504),@(0.450009765625,51.47611203514024),@(0.49952183894121072,54.393964482663613),@(0.54932484371947043,57.711395193162211),@(0.59987771594002193,61.343723298054684),@(0.6132573854161304,62.539558439217842),@(0.62683866021608314,63.793438110225999),@(0.64014058230669524,65.094393989811252),@(0.65281551932676074,66.378532583145684),@(0.66623064879083416,67.762470191850198),@(0.67855385122352663,69.156763250204108),@(0.69186186939045646,70.568007301422512),@(0.70480475512773422,72.008178716829818),@(0.71757888953986826,73.585128859709059),@(0.72921354495469559,75.019576425249994),@(0.7409401958050934,76.530757086478957),@(0.75228916849994132,78.269147021443175),@(0.76334982130392481,80.059175636229156),@(0.77419236031243888,81.883503928365769),@(0.7851260478131955,83.816580574711182),@(0.79605812346509619,85.708138105225331),@(0.80634007773235239,87.890532543024307),@(0.81664849632475134,90.057698662021082),@(0.8265550193224771,92.464762002472057),@(0.83636431184420357,94.839217638968705),@(0.84722632038814771,97.57794899358376),@(0.85666529344358933,100.24158747798405),@(0.86569413865232725,103.21319750833669),@(0.87416809033559761,106.40969503453539),@(0.88334306853582556,109.69455746672806),@(0.89171370882645395,113.28470368738286),@(0.90012886877265474,117.3147693592414),@(0.97360841729458814,201.07762495135049),@(0.97899290563827113,219.25214593472057),@(0.98763923688181277,270.70623602418016),@(0.99209903819206502,322.94745812676013),@(0.99265046495675524,331.09447803577751))
for ($i = 0; $i -lt $dataKerogenN2.Count; $i++) {
    $r = $i + 2
    $sheetKerogenN2.Cells.Item($r, 1).Value = $dataKerogenN2[$i][0]
    $sheetKerogenN2.Cells.Item($r, 2).Value = $dataKerogenN2[$i][1]
}

# --- Step 5: add "1_223_kerogen_Ar" sheet at the end (becomes the active sheet) ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetKerogenAr = $wb.Worksheets.Add($null, $lastSheet2)
$sheetKerogenAr.Name = "1_223_kerogen_Ar"
$sheetKerogenAr.Cells.Item(1,1).Value = "p_rels"
$sheetKerogenAr.Cells.Item(1,2).Value = "q_abs"
$dataKerogenAr = @(@(0.30035589377417815,52.571591696904839),@(0.35194522049085075,56.279304316433425),@(0.39969265288034533,59.739841141420207),@(0.45012415919319165,63.488239313717379),@(0.49947729770145866,67.393691305711002),@(0.54935387389567159,71.631260580204739),@(0.59937364916865066,76.43755426948168),@(0.64872423838742388,81.690630948754688),@(0.69797208521740228,87.975412368040367),@(0.74767131213207094,95.709974907273704),@(0.79700683760795576,105.70348102084671),@(0.82356379067037888,112.52051587474033),@(0.84858404237334506,120.56034663543336),@(0.87250918518891762,130.07385737456764),@(0.8965256616526005,142.39983834985941),@(0.92084885224755175,158.87843884429617),@(0.94857153338927436,189.13090605152126),@(0.95804263910360865,204.44244336656558),@(0.96661116187918417,224.4564731482709),@(0.97475835748665551,247.94420413181876),@(0.98423011514120828,291.91300217733379),@(0.98929717658768757,331.56311915239809),@(0.99120273467814202,349.61352603441293),@(0.9922860705472033,367.46853766804099))
for ($i = 0; $i -lt $dataKerogenAr.Count; $i++) {
    $r = $i + 2
    $sheetKerogenAr.Cells.Item($r, 1).Value = $dataKerogenAr[$i][0]
    $sheetKerogenAr.Cells.Item($r, 2).Value = $dataKerogenAr[$i][1]
}
$sheetKerogenAr.Activate()
$sheetKerogenAr.Range("O17").Select()
